# Batch import template rework: drop school-related columns, keep only
# the siswa (student) import fields, and rename "Jenis Kelamin" -> "Jenis
# Kelamin (L/P)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current header layout (A..P):
#  A NIP | B NPSN | C Nama Sekolah | D Alamat Sekolah | E Kabupaten Asal Sekolah
#  F Kecamatan Asal Sekolah | G Status (Swasta/Negeri) | H No. KK | I No. PKH
#  J No. PIP | K Nama Siswa | L Tempat Lahir | M Tanggal Lahir | N Alamat
#  O Jenis Kelamin | P Status (WNI/WNA)
#
# Target header layout (A..J):
#  A NIP | B No. KK | C No. PKH | D No. PIP | E Nama Siswa | F Tempat Lahir
#  G Tanggal Lahir | H Alamat | I Jenis Kelamin (L/P) | J Status (WNI/WNA)

# Rename the header text before the column moves, so the new text travels
# along with the rest of that column's data/formatting.
$ws.Range("O1").Value = "Jenis Kelamin (L/P)"

# Widen that column to fit the new, longer label (matches the "bestFit"
# recompute Excel performs automatically when a header no longer fits).
$ws.Columns.Item(15).ColumnWidth = 17.166666666666668

# Drop the contiguous block of school-related columns (B:G) -- NPSN, Nama
# Sekolah, Alamat Sekolah, Kabupaten Asal Sekolah, Kecamatan Asal Sekolah,
# Status (Swasta/Negeri). Everything to the right shifts left to fill in.
$ws.Range("B1:G1").EntireColumn.Delete()

$ws.Range("I2").Select()
